$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (new row): "Code Arduino" task ------------------------------
# Reuses the shared-string slot vacated by the old "moteur axe xy" text.
$ws.Range("A4").Value = "Code Arduino"

# --- Row 3: rename task to the fuller description ----------------------
$ws.Range("A3").Value = "Etude/Construction Moteurs Axe x/y"

# --- Cell fills ----------------------------------------------------------
# Red fill (FFFF0000) + red font on E4 -> new fill #4 / cellXfs #3
$ws.Range("E4").Interior.Color = 255
$ws.Range("E4").Font.Color = 255

# Same red fill, default font, on F4 -> reuses fill #4 / new cellXfs #4
$ws.Range("F4").Interior.Color = 255

# Green fill (FF00B050), default font, on E3 -> new fill #5 / cellXfs #5
$ws.Range("E3").Interior.Color = 5287936

# --- Final selection matches where the user ended up editing -----------
$ws.Range("A4").Select() | Out-Null
